$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.774.61'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').Value = '1.644.09'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  +0.44%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.63'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.499'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.81%  '
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('E8').Value = '  +0.33%  '
$ws.Range('E9').Value = '  -0.86%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.16'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.18%  '
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').Value = '1.869.20'
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').Value = '1.644.11'
$ws.Range('E13').Value = '  -0.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.16'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.528'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.30'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -3.03%  '
$ws.Range('D17').Value = '26.789.59'
$ws.Range('E17').Value = '  +0.25%  '
$ws.Range('E18').Value = '  -2.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '213.67'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.64%  '
$ws.Range('E20').Value = '  +0.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.36'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.28'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.37'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +11.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.36'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '145.03'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.52%  '
$ws.Range('E26').Value = '  +0.64%  '
$ws.Range('E27').Value = '  -1.95%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.08'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.65'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.70%  '
$ws.Range('E30').Value = '  -1.10%  '
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.31'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.85%  '
$ws.Range('E33').Value = '  -1.87%  '
$ws.Range('D34').Value = '1.287.39'
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  -1.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.43'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0173'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -5.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.537'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.88%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.826'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.32%  '
$ws.Range('E40').Value = '  +0.45%  '
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.35'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.64%  '
$ws.Range('D44').Value = '1.795.45'
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '61.10'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.38'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.38%  '
$ws.Range('E47').Value = '  -0.92%  '
$ws.Range('E48').Value = '  -1.49%  '
$ws.Range('E49').Value = '  +0.25%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.70'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.50%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0979'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.02%  '
